# g2.3 - ajuste no nome da coluna
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("g2.3")

# Rows 2-7: "Quantidade 2023/2010" -> "Quantidade (atual/dez anos antes)" with updated values
$ws.Range("B2:B7").Value = "Quantidade (atual/dez anos antes)"
$ws.Range("C2").Value = 69.07756813417191
$ws.Range("C3").Value = 51.02146541432979
$ws.Range("C4").Value = -19.79246840347153
$ws.Range("C5").Value = -61.70297632359825
$ws.Range("C6").Value = -70.15466408893185
$ws.Range("C7").Value = 110.3235747303544

# Rows 8-13: "Valor 2023/2010" -> "Valor (atual/dez anos antes)" with updated values
$ws.Range("B8:B13").Value = "Valor (atual/dez anos antes)"
$ws.Range("C8").Value = 39.20704066561322
$ws.Range("C9").Value = 111.3766687945756
$ws.Range("C10").Value = -38.80258191846784
$ws.Range("C11").Value = -75.00701801464335
$ws.Range("C12").Value = -66.55807726662741
$ws.Range("C13").Value = 33.87882696012225

# Rows 14-19: "Quantidade 2023/2022" -> "Quantidade (atual/ano anterior)" (text only)
$ws.Range("B14:B19").Value = "Quantidade (atual/ano anterior)"

# Rows 20-25: "Valor 2023/2022" -> "Valor (atual/ano anterior)" (text only)
$ws.Range("B20:B25").Value = "Valor (atual/ano anterior)"
